$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the match data (columns F:V) between row 36 and row 37 ---
# (columns A:E - Indice/pais/torneio/temporada/data_partida - stay put)
$row36 = @{}
$row37 = @{}
for ($c = 6; $c -le 22; $c++) {
    $row36[$c] = $ws.Cells.Item(36, $c).Value2
    $row37[$c] = $ws.Cells.Item(37, $c).Value2
}
for ($c = 6; $c -le 22; $c++) {
    $ws.Cells.Item(36, $c).Value = $row37[$c]
    $ws.Cells.Item(37, $c).Value = $row36[$c]
}

# --- 2. Append two new match rows (76 and 77) ---
# Copy formatting from the last existing row (75) so styles (bold index,
# date format, etc.) match the rest of the sheet.
$ws.Range("A75:V75").Copy()
$ws.Range("A76:V77").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{ row = 76; A = 75; B = "bosnia-and-herzegovina"; C = "premijer-liga-bih"; D = "2023-2024"; E = 45235.54166666666;
       F = "Siroki Brijeg"; G = 1; H = "Posusje"; I = 1; J = 1.9; K = "04/11/2023 01:13"; L = 2.44; M = "05/11/2023 12:57";
       N = 3.07; O = "04/11/2023 01:13"; P = 2.69; Q = "05/11/2023 12:57"; R = 3.8; S = "04/11/2023 01:13"; T = 3.44; U = "05/11/2023 12:57";
       V = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/siroki-brijeg-posusje/tpgw3CwO/" },
    @{ row = 77; A = 76; B = "bosnia-and-herzegovina"; C = "premijer-liga-bih"; D = "2023-2024"; E = 45235.54166666666;
       F = "Zrinjski"; G = 3; H = "Tuzla City"; I = 1; J = 1.21; K = "04/11/2023 01:13"; L = 1.16; M = "05/11/2023 12:44";
       N = 5.75; O = "04/11/2023 01:13"; P = 7.26; Q = "05/11/2023 12:55"; R = 8.65; S = "04/11/2023 01:13"; T = 14.64; U = "05/11/2023 12:55";
       V = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zrinjski-tuzla-city/YJkV4Y8B/" }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($rowData in $newRows) {
    $r = $rowData["row"]
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value = $rowData[$col]
    }
}

Write-Host "edit applied"
